# Adds two new columns (I = "I0", J = "IF") to the season-log sheet,
# mirroring the header style used by the existing "IP" column (H1) and
# filling in the per-game numeric values for rows 2-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting of the existing last header cell (H1: bold, bordered,
# centered) onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows -----------------------------------------------------------
# Each row: row number, I value, J value
$data = @(
    @(2,8,8),
    @(3,8,8),
    @(4,6,7),
    @(5,8,8),
    @(6,8,9),
    @(7,7,8),
    @(8,7,7),
    @(9,8,8),
    @(10,8,8),
    @(11,10,10),
    @(12,8,9),
    @(13,8,8),
    @(14,7,8),
    @(15,8,8),
    @(16,7,7),
    @(17,6,6),
    @(18,9,9),
    @(19,9,9),
    @(20,9,9),
    @(21,8,8),
    @(22,7,7),
    @(23,9,9),
    @(24,6,6),
    @(25,7,7),
    @(26,5,5),
    @(27,7,7),
    @(28,8,8),
    @(29,4,5),
    @(30,1,1),
    @(31,6,6),
    @(32,6,6),
    @(33,6,6),
    @(34,9,9),
    @(35,7,7),
    @(36,9,9),
    @(37,9,9),
    @(38,8,8),
    @(39,4,5),
    @(40,6,7),
    @(41,7,7),
    @(42,7,7),
    @(43,5,6),
    @(44,7,7),
    @(45,5,5),
    @(46,6,7),
    @(47,7,7),
    @(48,6,6),
    @(49,8,8),
    @(50,6,6),
    @(51,5,6),
    @(52,5,5),
    @(53,5,5),
    @(54,9,9),
    @(55,6,7),
    @(56,6,6),
    @(57,8,8),
    @(58,8,8),
    @(59,7,7),
    @(60,6,6),
    @(61,8,8),
    @(62,7,7),
    @(63,6,6),
    @(64,8,8),
    @(65,7,8),
    @(66,9,9),
    @(67,5,6),
    @(68,7,7),
    @(69,5,5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
